# Add "Behaviour change intervention style of delivery" (BCI style of delivery)
# to the upper level BCIO class-definitions sheet, and refresh a couple of
# related bits of content (column E header + "Intervention" row's example).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column header E1: "Elaboration" -> "Examples" -------------------------
$ws.Cells.Item(1, 5).Value = "Examples"

# --- Insert a new row 30 for the "BCI style of delivery" class -------------
# (this pushes the existing rows 30-42 down to 31-43, inheriting the A2 style)
$ws.Rows.Item(30).Insert()

$ws.Cells.Item(30, 1).Value = "BCIO:044000"
$ws.Cells.Item(30, 2).Value = "Behaviour change intervention style of delivery (BCI style of delivery)"
$ws.Cells.Item(30, 3).Value = "An attribute of a BCI delivery that encompasses the characteristics of how a BCI is communicated."
$ws.Cells.Item(30, 4).Value = "Process attribute/Process"
$ws.Cells.Item(30, 5).Value = "An example is cold and distant vs. warm and accepting."

# --- "Intervention" row (now row 37 after the insert): replace the old
# placeholder elaboration note with a proper "Examples" entry -------------
$ws.Cells.Item(37, 5).Value = "Examples of interventions are putting health warnings on cigarette packets, providing free stop smoking services and banning smoking in public places."

# --- View/selection tidy-up to mirror the saved workbook state -------------
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E2").Select()
